$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $NewValue)
    $r = $ws.Range($CellRef)
    $r.NumberFormat = "@"
    $r.Value2 = $NewValue
    $r.Style = "Normal"
}

$ws.Range('D2').Value2 = '34.640.59'
$ws.Range('E2').Value2 = '  +1.21%  '
$ws.Range('D3').Value2 = '1.801.27'
$ws.Range('E3').Value2 = '  +1.03%  '
$ws.Range('E4').Value2 = '  -0.13%  '
Set-TextValue 'D5' '227.55'
$ws.Range('E5').Value2 = '  +0.51%  '
Set-TextValue 'D6' '0.558'
$ws.Range('E6').Value2 = '  +1.92%  '
Set-TextValue 'D8' '32.89'
$ws.Range('E8').Value2 = '  +3.08%  '
$ws.Range('E9').Value2 = '  +1.61%  '
Set-TextValue 'D10' '0.0700'
$ws.Range('E10').Value2 = '  +1.37%  '
Set-TextValue 'D11' '0.0950'
$ws.Range('E11').Value2 = '  +0.32%  '
$ws.Range('D12').Value2 = '2.060.89'
$ws.Range('E12').Value2 = '  +1.05%  '
$ws.Range('B13').Value2 = 'Chainlink'
$ws.Range('C13').Value2 = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D13' '11.20'
$ws.Range('E13').Value2 = '  +1.37%  '
$ws.Range('B14').Value2 = 'WrappedEther'
$ws.Range('C14').Value2 = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value2 = '1.803.25'
$ws.Range('E14').Value2 = '  +1.02%  '
Set-TextValue 'D15' '0.639'
$ws.Range('E15').Value2 = '  +2.50%  '
$ws.Range('D16').Value2 = '34.606.88'
$ws.Range('E16').Value2 = '  +1.20%  '
$ws.Range('E17').Value2 = '  +3.08%  '
Set-TextValue 'D18' '69.00'
$ws.Range('E18').Value2 = '  +1.67%  '
$ws.Range('D19').Value2 = '0.0₃0808'
$ws.Range('E19').Value2 = '  +1.28%  '
Set-TextValue 'D20' '247.83'
$ws.Range('E20').Value2 = '  +0.31%  '
Set-TextValue 'D21' '11.32'
$ws.Range('E21').Value2 = '  +2.92%  '
$ws.Range('E22').Value2 = '  -0.13%  '
$ws.Range('E23').Value2 = '  +2.25%  '
Set-TextValue 'D24' '167.50'
$ws.Range('E24').Value2 = '  +3.20%  '
$ws.Range('E25').Value2 = '  +1.33%  '
$ws.Range('E26').Value2 = '  +1.58%  '
Set-TextValue 'D27' '16.62'
$ws.Range('E27').Value2 = '  +1.89%  '
Set-TextValue 'D28' '0.116'
$ws.Range('E28').Value2 = '  +2.31%  '
$ws.Range('E29').Value2 = '  -0.26%  '
$ws.Range('E30').Value2 = '  +11.29%  '
$ws.Range('B31').Value2 = 'PancakeSwap'
$ws.Range('C31').Value2 = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D31' '1.24'
$ws.Range('E31').Value2 = '  +0.87%  '
$ws.Range('B32').Value2 = 'Hedera'
$ws.Range('C32').Value2 = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D32' '0.0526'
$ws.Range('E32').Value2 = '  +1.04%  '
$ws.Range('E33').Value2 = '  +2.04%  '
$ws.Range('E34').Value2 = '  +2.66%  '
$ws.Range('D35').Value2 = '1.433.02'
$ws.Range('E35').Value2 = '  -0.78%  '
Set-TextValue 'D36' '2.60'
$ws.Range('E36').Value2 = '  +7.63%  '
Set-TextValue 'D37' '0.675'
Set-TextValue 'D38' '1.08'
$ws.Range('E38').Value2 = '  +3.19%  '
$ws.Range('E39').Value2 = '  +0.58%  '
Set-TextValue 'D40' '85.57'
$ws.Range('E40').Value2 = '  +6.51%  '
$ws.Range('E41').Value2 = '  +1.37%  '
Set-TextValue 'D42' '0.942'
$ws.Range('E42').Value2 = '  +1.91%  '
Set-TextValue 'D43' '2.77'
$ws.Range('E43').Value2 = '  +3.91%  '
Set-TextValue 'D44' '13.76'
$ws.Range('E44').Value2 = '  +1.01%  '
$ws.Range('E45').Value2 = '  +3.50%  '
$ws.Range('E46').Value2 = '  +0.59%  '
$ws.Range('E47').Value2 = '  +0.34%  '
$ws.Range('D48').Value2 = '1.959.34'
$ws.Range('E48').Value2 = '  +0.94%  '
Set-TextValue 'D49' '106.26'
$ws.Range('E49').Value2 = '  +1.37%  '
$ws.Range('E50').Value2 = '  -0.15%  '
$ws.Range('D51').Value2 = '0.0₆0129'
$ws.Range('E51').Value2 = '  -5.67%  '
